# Update "想去人数" (interested-count) values in the F column across all four sheets.
# Generated from the commit diff (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 43
$ws.Range("F4").Value = 3596
$ws.Range("F5").Value = 3596
$ws.Range("F7").Value = 5129
$ws.Range("F8").Value = 536
$ws.Range("F9").Value = 365
$ws.Range("F11").Value = 695
$ws.Range("F13").Value = 97
$ws.Range("F15").Value = 705
$ws.Range("F16").Value = 320
$ws.Range("F18").Value = 93
$ws.Range("F19").Value = 157
$ws.Range("F22").Value = 4925
$ws.Range("F26").Value = 6057
$ws.Range("F28").Value = 17
$ws.Range("F29").Value = 3225
$ws.Range("F30").Value = 347
$ws.Range("F31").Value = 715
$ws.Range("F32").Value = 4444
$ws.Range("F36").Value = 1034
$ws.Range("F40").Value = 875
$ws.Range("F41").Value = 1014
$ws.Range("F42").Value = 2031

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 13
$ws.Range("F5").Value = 55

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 229
$ws.Range("F3").Value = 1123
$ws.Range("F4").Value = 50

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 229
$ws.Range("F3").Value = 43
$ws.Range("F4").Value = 1123
$ws.Range("F5").Value = 50
$ws.Range("F7").Value = 3596
$ws.Range("F8").Value = 3596
$ws.Range("F10").Value = 5129
$ws.Range("F11").Value = 536
$ws.Range("F12").Value = 365
$ws.Range("F14").Value = 695
$ws.Range("F16").Value = 97
$ws.Range("F18").Value = 705
$ws.Range("F19").Value = 320
$ws.Range("F22").Value = 93
$ws.Range("F23").Value = 157
$ws.Range("F26").Value = 4925
$ws.Range("F30").Value = 6057
$ws.Range("F32").Value = 17
$ws.Range("F33").Value = 3225
$ws.Range("F34").Value = 347
$ws.Range("F35").Value = 715
$ws.Range("F36").Value = 4444
$ws.Range("F41").Value = 1034
$ws.Range("F45").Value = 875
$ws.Range("F46").Value = 1014
$ws.Range("F47").Value = 13
$ws.Range("F48").Value = 2031
$ws.Range("F50").Value = 55

